$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45987
$ws.Range("B2").Value = 78.14
$ws.Range("C2").Value = 68.97
$ws.Range("D2").Value = 62.38
$ws.Range("E2").Value = 63.08
$ws.Range("F2").Value = 61.99
$ws.Range("G2").Value = 73.91
$ws.Range("H2").Value = 81.05
$ws.Range("I2").Value = 94.54000000000001
$ws.Range("J2").Value = 104.67
$ws.Range("K2").Value = 69.2
$ws.Range("L2").Value = 33.79
$ws.Range("M2").Value = 19.55
$ws.Range("N2").Value = 29.77
$ws.Range("O2").Value = 29.12
$ws.Range("P2").Value = 29.84
$ws.Range("Q2").Value = 32.49
$ws.Range("R2").Value = 69.95999999999999
$ws.Range("S2").Value = 97.45999999999999
$ws.Range("T2").Value = 120.43
$ws.Range("U2").Value = 142.01
$ws.Range("V2").Value = 124.83
$ws.Range("W2").Value = 100.21
$ws.Range("X2").Value = 91.89
$ws.Range("Y2").Value = 88.89
$ws.Range("Z2").Value = 73.67
$ws.Range("AA2").Value = "16h-20h"
$ws.Range("AB2").Value = 107.46
$ws.Range("AD2").Value = 131.22
$ws.Range("AF2").Value = 112.52
$ws.Range("AG2").Value = "1h-16h"
